$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.3921459680713356
    "C2" = 0.4690580135690283
    "B3" = 0.1938669551675961
    "C3" = 0.1186057994486912
    "B4" = 0.0677949822057861
    "C4" = 0.06717876075738202
    "B5" = 0.06719666705145336
    "C5" = 0.0671787607573822
    "B6" = 0.07195801386525061
    "C6" = 0.06717876075738212
    "B7" = 0.08419062385962513
    "C7" = 0.08023122541279935
    "B8" = 0.1228467897789531
    "C8" = 0.130568679297338
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
